$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set all affected D/E cells to text format so values like "1.00" are not coerced to numbers
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.534.48"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -3.62%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.997.40"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -2.62%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "537.97"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.28%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "132.91"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -1.94%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "2.991.31"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -2.67%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.493"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.24%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.147"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -5.33%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.09"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -6.44%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.447"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.68%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000221"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -2.37%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.76"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.63%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.483.43"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -2.61%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "61.604.89"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -3.60%  "
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -2.24%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.998.68"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -2.72%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.61"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.52%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "469.39"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -2.73%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.15"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -1.45%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.669"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -4.16%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.97"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -1.55%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "80.17"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.55%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.94"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.61%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.998"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.15%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.69"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.20%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.70"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -4.96%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.05%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.89"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.70%  "
$ws.Range("B31").NumberFormat = "@"
$ws.Range("B31").Value = "Mantle"
$ws.Range("C31").NumberFormat = "@"
$ws.Range("C31").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.15"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -2.18%  "
$ws.Range("B32").NumberFormat = "@"
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").NumberFormat = "@"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "25.52"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -2.68%  "
$ws.Range("B33").NumberFormat = "@"
$ws.Range("B33").Value = "OKB"
$ws.Range("C33").NumberFormat = "@"
$ws.Range("C33").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "55.49"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -2.21%  "
$ws.Range("B34").NumberFormat = "@"
$ws.Range("B34").Value = "Stacks"
$ws.Range("C34").NumberFormat = "@"
$ws.Range("C34").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.29"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -2.67%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.36"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.48%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.89"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.64%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "455.35"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -8.62%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.181.40"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.74%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0789"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.03%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.119"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.36%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0380"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -4.26%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.08"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.05%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.38"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -9.84%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.05%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "25.53"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +5.05%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.242"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -4.65%  "
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = "Monero"
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "118.37"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -2.86%  "
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = "Fetch.AI"
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.97"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -3.16%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.107"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.96%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0$([char]0x2083)0487"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -7.67%  "
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = "BitgetToken"
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = "https://coinranking.com/coin/q7gMmMdLb+bitgettoken-bgb"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.25"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +5.62%  "